$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D3: fix "Eigen-," -> "Eigenbewirtschaftung," inside the Erwerbsart/Flaechenaufstellung combo string
$ws.Range("D3").Value = "Erwerbsart:select(Vollerwerb, Nebenerwerb):pflicht;Flächenaufstellung:select(Eigenbewirtschaftung,Pacht,Mitbewirtschaftung):pflicht"

# C7: fix typo "Nebetätigkeiten" -> "Nebentätigkeiten"
$ws.Range("C7").Value = "A3. Nebentätigkeiten gesamt (Erstaufnahme)"

# Remove the stray ParentID values in B29 / B30 (cells fully cleared, not just blanked)
$ws.Range("B29").Clear()
$ws.Range("B30").Clear()

# D32: add new "Bemerkung:text" attribute cell
$ws.Range("D32").Value = "Bemerkung:text"

# F34: add new "Ja" upload flag cell, matching the text-formatted style used throughout column F
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = "Ja"
